$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.258.85'
$ws.Range('E2').Value = '  +1.28%  '

# Row 3
$ws.Range('D3').Value = '1.655.34'
$ws.Range('E3').Value = '  +1.07%  '

# Row 4
$ws.Range('E4').Value = '  +1.31%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.29'
$ws.Range('E5').Value = '  +0.86%  '

# Row 6
$ws.Range('E6').Value = '  +0.95%  '

# Row 7
$ws.Range('E7').Value = '  +1.34%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0640'
$ws.Range('E8').Value = '  +0.28%  '

# Row 9
$ws.Range('E9').Value = '  -0.48%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.49'
$ws.Range('E10').Value = '  -0.60%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0801'
$ws.Range('E11').Value = '  +0.89%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.716.97'
$ws.Range('E12').Value = '  +4.75%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.28'
$ws.Range('E13').Value = '  +0.35%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.544'
$ws.Range('E14').Value = '  -0.03%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.64'
$ws.Range('E15').Value = '  +1.01%  '

# Row 16
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  -0.12%  '

# Row 17
$ws.Range('D17').Value = '26.212.92'
$ws.Range('E17').Value = '  +1.02%  '

# Row 18
$ws.Range('E18').Value = '  +1.32%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '194.52'
$ws.Range('E19').Value = '  +0.86%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.34'
$ws.Range('E20').Value = '  -0.45%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.79'
$ws.Range('E21').Value = '  -1.34%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.20'
$ws.Range('E22').Value = '  -1.35%  '

# Row 23
$ws.Range('E23').Value = '  +1.54%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.14'
$ws.Range('E24').Value = '  +0.67%  '

# Row 25
$ws.Range('E25').Value = '  +1.71%  '

# Row 26
$ws.Range('E26').Value = '  -0.36%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.88'
$ws.Range('E27').Value = '  +0.58%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.57'
$ws.Range('E28').Value = '  +0.02%  '

# Row 29
$ws.Range('E29').Value = '  +0.68%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0490'
$ws.Range('E30').Value = '  -2.59%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.27'
$ws.Range('E31').Value = '  +1.08%  '

# Row 32
$ws.Range('E32').Value = '  -0.91%  '

# Row 33
$ws.Range('E33').Value = '  +0.17%  '

# Row 34
$ws.Range('E34').Value = '  +1.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.906'
$ws.Range('E35').Value = '  +0.50%  '

# Row 36
$ws.Range('D36').Value = '1.138.95'
$ws.Range('E36').Value = '  +0.15%  '

# Row 37
$ws.Range('E37').Value = '  +0.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.533'
$ws.Range('E38').Value = '  -2.28%  '

# Row 39
$ws.Range('E39').Value = '  -0.46%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.803'
$ws.Range('E40').Value = '  +0.64%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.93'
$ws.Range('E41').Value = '  -0.37%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.34'
$ws.Range('E42').Value = '  -2.57%  '

# Row 43
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').Value = '0.0₆0113'
$ws.Range('E43').Value = '  -0.95%  '

# Row 44
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '56.61'
$ws.Range('E44').Value = '  -0.07%  '

# Row 45
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.50'
$ws.Range('E45').Value = '  +1.75%  '

# Row 46
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0524'
$ws.Range('E46').Value = '  -1.41%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.78'
$ws.Range('E47').Value = '  +1.13%  '

# Row 48
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.419'
$ws.Range('E48').Value = '  +1.07%  '

# Row 49
$ws.Range('B49').Value = 'USDD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.01'
$ws.Range('E49').Value = '  +1.17%  '

# Row 50
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0941'
$ws.Range('E50').Value = '  -2.38%  '

# Row 51
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').Value = '  +2.13%  '
